# Update "simulated-data" workbook per commit "update data with new whatdataio version"
#
# Changes:
#  - Rebrand "NCC Priority Actions App" -> "What To Do application" in the two
#    worksheet description cells that mention it (Site data!B2, Feasibility data!B2).
#  - Site data: Lobster Bay & Round Bay current status changes from "Restore" to
#    "Signage"; cost figures (columns E:G) refreshed for all four sites.
#  - Feature data: goal/weight figures (columns B:C) refreshed for all three features.
#  - Consequence of "Maintain" / "Signage" / "Restore": expectation figures
#    (columns B:D) refreshed for all four sites.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Unprotect()
}

# Worksheet order (fixed, from workbook.xml):
#  1 Site data | 2 Feasibility data | 3 Feature data
#  4 Consequence of "Maintain" | 5 Consequence of "Signage" | 6 Consequence of "Restore"
#  7 metadata
$siteData = $wb.Worksheets.Item(1)
$feasibilityData = $wb.Worksheets.Item(2)
$featureData = $wb.Worksheets.Item(3)
$maintain = $wb.Worksheets.Item(4)
$signage = $wb.Worksheets.Item(5)
$restore = $wb.Worksheets.Item(6)

# --- Site data -----------------------------------------------------------
$siteData.Range("B2").Value = 'Specifically, we ask that you input the longitude and latitude (in decimal degrees) of each site. If you have a shapefile with spatial locations (e.g. point localities, boundaries) of your sites, these can also be supplied in the What To Do application. We also ask that you input the cost of implementing each management action (e.g. in Canadian Dollars) within each site. Please note that cost values should not be below zero (though they can equal zero) and not be greater than 1,000,000 (i.e. one million). As such, you might need to rescale your cost values. For example, if one of your cost values is “10000000” Canadian Dollars, instead of inputting values as Canadian Dollars, you could you input values as thousands of Canadian Dollars (i.e. “10000”). Please take care to ensure that all cost values are in the same units. After filling out this worksheet, every light gray cell should have a numerical value.'

$siteData.Range("D5").Value = "Signage"
$siteData.Range("D7").Value = "Signage"

$siteData.Range("E4").Value = 174
$siteData.Range("F4").Value = 371
$siteData.Range("G4").Value = 596

$siteData.Range("E5").Value = 154
$siteData.Range("F5").Value = 296
$siteData.Range("G5").Value = 527

$siteData.Range("E6").Value = 147
$siteData.Range("F6").Value = 207
$siteData.Range("G6").Value = 363

$siteData.Range("E7").Value = 198
$siteData.Range("F7").Value = 364
$siteData.Range("G7").Value = 520

# --- Feasibility data ------------------------------------------------------
$feasibilityData.Range("B2").Value = 'Specifically, we ask that you input data indicating which management actions are feasible to implement within each site. By default, all actions can be potentially implemented within each and every site. To specify that a certain a certain action cannot be implemented within a certain site, please enter value of “0”. You can also ensure that a certain site can ONLY have a certain action implemented within it, by specifying a value of “0” for every other action. This information, if you prefer, can also be specified within the What To Do application---however, you will have to re-specify this information each and every time you open the application.'

# --- Feature data ----------------------------------------------------------
$featureData.Range("B4").Value = 95
$featureData.Range("C4").Value = 24

$featureData.Range("B5").Value = 24
$featureData.Range("C5").Value = 7

$featureData.Range("B6").Value = 27
$featureData.Range("C6").Value = 27

# --- Consequence of "Maintain" ----------------------------------------------
$maintain.Range("B4").Value = 155
$maintain.Range("C4").Value = 173
$maintain.Range("D4").Value = 110

$maintain.Range("B5").Value = 180
$maintain.Range("C5").Value = 163
$maintain.Range("D5").Value = 144

$maintain.Range("B6").Value = 181
$maintain.Range("C6").Value = 141
$maintain.Range("D6").Value = 157

$maintain.Range("B7").Value = 149
$maintain.Range("C7").Value = 170
$maintain.Range("D7").Value = 162

# --- Consequence of "Signage" -----------------------------------------------
$signage.Range("B4").Value = 206
$signage.Range("C4").Value = 206
$signage.Range("D4").Value = 291

$signage.Range("B5").Value = 379
$signage.Range("C5").Value = 207
$signage.Range("D5").Value = 237

$signage.Range("B6").Value = 396
$signage.Range("C6").Value = 298
$signage.Range("D6").Value = 267

$signage.Range("B7").Value = 287
$signage.Range("C7").Value = 229
$signage.Range("D7").Value = 251

# --- Consequence of "Restore" -----------------------------------------------
$restore.Range("B4").Value = 373
$restore.Range("C4").Value = 441
$restore.Range("D4").Value = 463

$restore.Range("B5").Value = 397
$restore.Range("C5").Value = 369
$restore.Range("D5").Value = 551

$restore.Range("B6").Value = 542
$restore.Range("C6").Value = 506
$restore.Range("D6").Value = 514

$restore.Range("B7").Value = 600
$restore.Range("C7").Value = 481
$restore.Range("D7").Value = 482

# Re-protect sheets to restore the original workbook's protection state.
foreach ($ws in $wb.Worksheets) {
    $ws.Protect()
}
